$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        B = 0.00228310502283105; C = 0.0045662100456621; D = 0.00228310502283105;
        E = 0.00228310502283105; F = 0.045662100456621; G = 0; H = 0;
        I = 0.0045662100456621; J = 0.0045662100456621; K = 0.0114155251141553;
        L = 0.00684931506849315; M = 0.0045662100456621; N = 0;
        O = 0.0182648401826484; P = 0.0045662100456621; Q = 0.00684931506849315;
        R = 0.988584474885845; S = 0; T = 0.0045662100456621; U = 0;
        V = 0.0045662100456621; W = 0.00684931506849315; X = 0.0136986301369863
    }
    3 = @{
        B = 0.00684931506849315; C = 0.972602739726027; D = 0.977168949771689;
        E = 0.984018264840183; F = 0.0045662100456621; G = 0.0159817351598174;
        H = 0.974885844748858; I = 0.00228310502283105; J = 0.0045662100456621;
        K = 0.0045662100456621; L = 0; M = 0.0091324200913242;
        N = 0.00228310502283105; O = 0; P = 0; Q = 0;
        R = 0.00684931506849315; S = 1; T = 0.0114155251141553;
        U = 0.974885844748858; V = 0; W = 0.00684931506849315;
        X = 0.00228310502283105
    }
    4 = @{
        B = 0.970319634703196; C = 0; D = 0.0159817351598174; E = 0.0045662100456621;
        F = 0.949771689497717; G = 0.00228310502283105; H = 0.0114155251141553;
        I = 0.00684931506849315; J = 0.988584474885845; K = 0.0045662100456621;
        L = 0.00684931506849315; M = 0.974885844748858; N = 0.0091324200913242;
        O = 0.970319634703196; P = 0.0091324200913242; Q = 0.990867579908676;
        R = 0; S = 0; T = 0; U = 0.00684931506849315; V = 0.0136986301369863;
        W = 0.986301369863014; X = 0.977168949771689
    }
    5 = @{
        B = 0.0114155251141553; C = 0.0228310502283105; D = 0.0045662100456621;
        E = 0.0091324200913242; F = 0; G = 0.981735159817352; H = 0.0136986301369863;
        I = 0.986301369863014; J = 0.00228310502283105; K = 0.979452054794521;
        L = 0.986301369863014; M = 0.0114155251141553; N = 0.988584474885845;
        O = 0.0114155251141553; P = 0.986301369863014; Q = 0.00228310502283105;
        R = 0.0045662100456621; S = 0; T = 0.984018264840183;
        U = 0.00684931506849315; V = 0.981735159817352; W = 0;
        X = 0.00684931506849315
    }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
